$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to stay text so values like "228.44" or
# "3.60" are not silently re-interpreted as numbers (which would drop
# trailing zeros / change formatting). Column D is the only one at risk;
# B, C and E values never round-trip as plain numbers.
foreach ($addr in @("D2", "D3", "D5", "D7", "D10", "D13", "D14", "D15", "D16", "D17", "D18", "D22", "D26", "D28", "D31", "D34", "D35", "D36", "D38", "D41", "D42", "D43", "D44", "D48", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "38.720.15"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "2.100.49"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "228.44"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "62.08"
$ws.Range("E7").Value = "  +1.56%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +2.08%  "
$ws.Range("D10").Value = "0.0843"
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("E12").Value = "  +7.10%  "
$ws.Range("D13").Value = "2.411.79"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "22.09"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "0.804"
$ws.Range("E15").Value = "  +3.48%  "
$ws.Range("D16").Value = "5.53"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "2.101.49"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "38.717.73"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").Value = "227.49"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").Value = "172.26"
$ws.Range("E26").Value = "  +1.56%  "
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").Value = "0.139"
$ws.Range("E28").Value = "  +6.63%  "
$ws.Range("E29").Value = "  +4.13%  "
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").Value = "2.53"
$ws.Range("E31").Value = "  +7.14%  "
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("D34").Value = "4.75"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("B35").Value = "THORChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D35").Value = "6.92"
$ws.Range("E35").Value = "  +7.60%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.0618"
$ws.Range("E36").Value = "  +2.57%  "
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("D38").Value = "3.60"
$ws.Range("E38").Value = "  +3.25%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("D41").Value = "102.69"
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("D42").Value = "0.0229"
$ws.Range("E42").Value = "  +4.49%  "
$ws.Range("D43").Value = "1.533.89"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "7.84"
$ws.Range("E44").Value = "  +4.66%  "
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("E46").Value = "  +3.24%  "
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").Value = "4.13"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").Value = "2.292.78"
$ws.Range("E51").Value = "  +0.05%  "
